# Applies the 2023-04-03 GitHub Actions cryptos-list refresh to Sheet1.
# D = Price, E = Volume(1h); rows 39/40 additionally swap Coin/Link (ranking change).
# D/E columns store plain text (e.g. "1.003", "  -2.53%  "), so force text
# number format before assigning values to stop COM auto-coercing them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Row 39 / Row 40: ranking swapped between Aptos and InternetComputer(DFINITY) ---
Set-TextValue "B39" "InternetComputer(DFINITY)"
Set-TextValue "C39" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D39" "5.023"
Set-TextValue "E39" "  -3.49%  "

Set-TextValue "B40" "Aptos"
Set-TextValue "C40" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D40" "11.51"
Set-TextValue "E40" "  +2.20%  "

# --- Remaining Price / Volume(1h) refresh across the rest of the list ---
Set-TextValue "D2" "27.810.43"
Set-TextValue "E2" "  -2.53%  "
Set-TextValue "D3" "1.782.08"
Set-TextValue "E3" "  -2.13%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "310.39"
Set-TextValue "E5" "  -1.92%  "
Set-TextValue "E6" "  -0.06%  "
Set-TextValue "D7" "0.5133"
Set-TextValue "E7" "  -0.11%  "
Set-TextValue "D8" "0.3793"
Set-TextValue "E8" "  -2.16%  "
Set-TextValue "D9" "0.07756"
Set-TextValue "E9" "  -8.21%  "
Set-TextValue "D10" "41.14"
Set-TextValue "E10" "  -1.74%  "
Set-TextValue "D11" "1.085"
Set-TextValue "E11" "  -2.36%  "
Set-TextValue "E12" "  +0.02%  "
Set-TextValue "D13" "6.199"
Set-TextValue "E13" "  -3.29%  "
Set-TextValue "E14" "  -4.73%  "
Set-TextValue "D15" "1.769.58"
Set-TextValue "E15" "  -2.38%  "
Set-TextValue "E16" "  -4.73%  "
Set-TextValue "D17" "91.40"
Set-TextValue "E17" "  -1.54%  "
Set-TextValue "D18" "0.00001071"
Set-TextValue "E18" "  -5.89%  "
Set-TextValue "D19" "0.06527"
Set-TextValue "E19" "  -2.50%  "
Set-TextValue "E20" "  +0.02%  "
Set-TextValue "D21" "17.01"
Set-TextValue "E21" "  -4.42%  "
Set-TextValue "D22" "5.910"
Set-TextValue "E22" "  -3.14%  "
Set-TextValue "D23" "27.861.24"
Set-TextValue "E23" "  -2.45%  "
Set-TextValue "E24" "  -4.07%  "
Set-TextValue "D25" "2.236"
Set-TextValue "E25" "  -1.86%  "
Set-TextValue "D26" "158.67"
Set-TextValue "E26" "  -0.21%  "
Set-TextValue "D27" "20.19"
Set-TextValue "E27" "  -4.43%  "
Set-TextValue "D28" "1.980.78"
Set-TextValue "E28" "  -2.16%  "
Set-TextValue "E29" "  -2.82%  "
Set-TextValue "D30" "125.37"
Set-TextValue "E30" "  -0.58%  "
Set-TextValue "E31" "  -0.47%  "
Set-TextValue "D32" "1.026"
Set-TextValue "E32" "  -6.50%  "
Set-TextValue "D33" "3.607"
Set-TextValue "E33" "  -1.83%  "
Set-TextValue "D34" "5.478"
Set-TextValue "E34" "  -4.86%  "
Set-TextValue "D35" "0.07090"
Set-TextValue "E35" "  -6.30%  "
Set-TextValue "D36" "0.02308"
Set-TextValue "E36" "  -2.52%  "
Set-TextValue "D37" "0.2121"
Set-TextValue "E37" "  -4.86%  "
Set-TextValue "E38" "  -1.05%  "
Set-TextValue "D41" "0.6084"
Set-TextValue "E41" "  -4.04%  "
Set-TextValue "E42" "  -0.09%  "
Set-TextValue "D43" "1.151"
Set-TextValue "E43" "  -3.57%  "
Set-TextValue "D44" "1.319"
Set-TextValue "E44" "  -5.84%  "
Set-TextValue "D45" "0.5960"
Set-TextValue "E45" "  +0.37%  "
Set-TextValue "D46" "13.06"
Set-TextValue "E46" "  -3.85%  "
Set-TextValue "D47" "3.709"
Set-TextValue "E47" "  -1.74%  "
Set-TextValue "D48" "126.89"
Set-TextValue "E48" "  +0.70%  "
Set-TextValue "E49" "  +1.21%  "
Set-TextValue "D50" "1.896"
Set-TextValue "E50" "  -4.87%  "
Set-TextValue "D51" "0.06708"
Set-TextValue "E51" "  -4.00%  "

